# Apply the "sample_inputs.xlsx" update:
#  - Aircraft speeds bumped from 60 -> 170 mph for both types
#  - Operator sheet: aircraft counts doubled/rebalanced (2/1 -> 8/8),
#    an extra charger column (H) wired up as "Yes"/"Slow", the Yes/No
#    COUNTIF range widened, and the one remaining "No" flag flipped to "Yes"
#  - Ports sheet: Landing Slots bumped from 1 -> 5 and a new Port 6 row added
#  - Ports sheet: drop the stray empty formatted cell I1
#  - Data validation range on Operator!C14:G14 grows to C14:H14
#  - Ports tab becomes the active tab/sheet

$wb = $excel.ActiveWorkbook

$operator = $wb.Worksheets.Item("Operator")
$aircraft = $wb.Worksheets.Item("Aircraft")
$ports    = $wb.Worksheets.Item("Ports")

# ---- Aircraft sheet: Speed (mph) 60 -> 170 for both aircraft types ----
$aircraft.Range("B2").Value = 170
$aircraft.Range("B3").Value = 170

# ---- Operator sheet: Number of Aircraft split 2/1 -> 8/8 ----
$operator.Range("C10").Value = 8
$operator.Range("D10").Value = 8

# ---- Operator sheet: Serviced Ports row - G13 flips from "No" to "Yes",
#      and a new H13 = "Yes" / H14 = "Slow" pair is added ----
$operator.Range("G13").Value = "Yes"

$operator.Range("H13").Value = "Yes"
$operator.Range("G13").Copy()
$operator.Range("H13").PasteSpecial(-4122)  # xlPasteFormats

$operator.Range("H14").Value = "Slow"
$operator.Range("G14").Copy()
$operator.Range("H14").PasteSpecial(-4122)  # xlPasteFormats

# Widen the Yes-count formula to include the new column
$operator.Range("B13").Formula = "=COUNTIF(C13:I13,""Yes"")"

# ---- Data validation: Slow/Fast/None list grows to cover column H ----
$slowFastRange = $operator.Range("C14:H14")
$slowFastRange.Validation.Delete()
$slowFastRange.Validation.Add(3, 1, 1, '"Slow, Fast, None"')

# ---- Ports sheet: Landing Slots 1 -> 5 for existing ports ----
$ports.Range("E2").Value = 5
$ports.Range("E3").Value = 5
$ports.Range("E4").Value = 5
$ports.Range("E5").Value = 5
$ports.Range("E6").Value = 5

# Drop the stray empty formatted cell at I1 (no longer used)
$ports.Range("I1").Clear()

# New Port 6 row
$ports.Range("A7").Value = 6
$ports.Range("B7").Value = 50
$ports.Range("C7").Value = 72
$ports.Range("D7").Value = 0.6
$ports.Range("E7").Value = 5
$ports.Range("F7").Value = 20
$ports.Range("D7").NumberFormat = $ports.Range("D6").NumberFormat
$ports.Range("F7").NumberFormat = $ports.Range("F6").NumberFormat

# ---- Make "Ports" the active sheet/tab ----
$operator.Range("A1").Select()
$aircraft.Range("A1").Select()
$ports.Range("A1").Select()
$ports.Activate()
